# Wijzigingen.docx - opgeruimd, AEX index B002 opgelost, aan sprint 038 begonnen
#
# Row "W002" (Tijdelijke dagkoers voor huidige handelsdag ophalen en tonen):
#   - status cell: add a remark about the intraday bug (B012)
#   - prio cell:   Hoog -> Matig

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 3 (1-based, header is row 1) is the W002 row.
# Column 4 = "status", column 5 = "prio".
$statusCell = $table.Cell(3, 4)
$statusRange = $statusCell.Range
$statusRange.Collapse(0)
$cr = [char]13
$statusRange.InsertAfter($cr + "Werkt niet intraday, zie fout B012")

$prioCell = $table.Cell(3, 5)
$prioCell.Range.Text = "Matig"
